$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nominatifs")

$lastRow = 21

# Insert a new column before column A: shifts Nom/Prenom data (and their
# column formatting, e.g. the bestFit width) from A:B to B:C.
$ws.Columns.Item(1).Insert()

# Fill the new column A with sequential numbers 1..20 for rows 2..21
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Add a new column D with role labels
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "acteur"
}
for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = "réalisateur"
}
for ($r = 12; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "client"
}

$ws.Range("D22").Select()
$ws.Activate()
